$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update row labels: drop the "%" suffix from the metric names
$ws.Range("A9").Value = "Acurácia"
$ws.Range("A10").Value = "Sensibilidade"
$ws.Range("A11").Value = "Especificidade"

$ws.Range("I9").Value = "Acurácia"
$ws.Range("I10").Value = "Sensibilidade"
$ws.Range("I11").Value = "Especificidade"

# Update the metric values: from percentage scale (100) to fraction scale (1)
$ws.Range("B9:E9").Value = 1
$ws.Range("B10:E10").Value = 1
$ws.Range("B11:E11").Value = 1

$ws.Range("J9:M9").Value = 1
$ws.Range("J10:M10").Value = 1
$ws.Range("J11:M11").Value = 1
